$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "All.jamais.sex"
$ws.Range("C1").Value = "Males.jamais.sex"
$ws.Range("D1").Value = "Females.jamais.sex"
$ws.Range("E1").Value = "Not known / missing.jamais.sex"
